$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 2900.3333
$ws.Range("I86").Value = 2800.7144
$ws.Range("J86").Value = 3039.8
$ws.Range("K86").Value = 2800.7144
$ws.Range("L86").Value = 3039.8
$ws.Range("M86").Value = -1677.7144
$ws.Range("N86").Value = -5285.8

# Row 89
$ws.Range("H89").Value = 2900.3333
$ws.Range("I89").Value = 2800.7144
$ws.Range("J89").Value = 3039.8
$ws.Range("K89").Value = 14003.572
$ws.Range("L89").Value = 15199
$ws.Range("M89").Value = -8387.572
$ws.Range("N89").Value = -26431

# Row 106
$ws.Range("H106").Value = 7939308.5
$ws.Range("I106").Value = 8549717
$ws.Range("K106").Value = 8549717
$ws.Range("M106").Value = -8549086

# Row 112
$ws.Range("H112").Value = 7143843.5
$ws.Range("J112").Value = 7576749.5
$ws.Range("L112").Value = 22730248.5
$ws.Range("N112").Value = -22732464.5

# Row 129
$ws.Range("H129").Value = 1145.3334
$ws.Range("J129").Value = 1249.5652
$ws.Range("L129").Value = 3748.6956
$ws.Range("N129").Value = -13748.6956

# Row 133
$ws.Range("H133").Value = 13516.154
$ws.Range("J133").Value = 13516.154
$ws.Range("L133").Value = 13516.154
$ws.Range("N133").Value = -23636.154

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 57074
$ws.Range("I2").Value = 63839.5
$ws.Range("K2").Value = 63839.5
$ws.Range("M2").Value = -63726.5

# Row 32
$ws.Range("H32").Value = 15902.04
$ws.Range("I32").Value = 2688.1077
$ws.Range("K32").Value = 2688.1077
$ws.Range("M32").Value = -2401.1077

# Row 61
$ws.Range("H61").Value = 1813.0513
$ws.Range("I61").Value = 1380.2354
$ws.Range("J61").Value = 4756.2
$ws.Range("K61").Value = 1380.2354
$ws.Range("L61").Value = 4756.2
$ws.Range("M61").Value = -1168.2354
$ws.Range("N61").Value = -5180.2

# Row 102
$ws.Range("H102").Value = 1808.8889

# Row 116
$ws.Range("H116").Value = 57074
$ws.Range("I116").Value = 63839.5
$ws.Range("K116").Value = 63839.5
$ws.Range("M116").Value = -61545.5

# Row 122
$ws.Range("H122").Value = 2264.182
$ws.Range("I122").Value = 2008.4
$ws.Range("J122").Value = 2477.3333
$ws.Range("K122").Value = 6025.200000000001
$ws.Range("L122").Value = 7431.999899999999
$ws.Range("M122").Value = -3575.200000000001
$ws.Range("N122").Value = -12331.9999

# Row 136
$ws.Range("H136").Value = 1813.0513
$ws.Range("I136").Value = 1380.2354
$ws.Range("J136").Value = 4756.2
$ws.Range("K136").Value = 4140.706200000001
$ws.Range("L136").Value = 14268.6
$ws.Range("M136").Value = -1590.706200000001
$ws.Range("N136").Value = -19368.6

# Row 139
$ws.Range("H139").Value = 42692.25
$ws.Range("J139").Value = 42692.25
$ws.Range("L139").Value = 42692.25
$ws.Range("N139").Value = -52972.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 57074
$ws.Range("I3").Value = 63839.5
$ws.Range("K3").Value = 63839.5
$ws.Range("M3").Value = -63725.5

# Row 5
$ws.Range("H5").Value = 15050.714
$ws.Range("I5").Value = 17816.666
$ws.Range("K5").Value = 17816.666
$ws.Range("M5").Value = -17703.666

# Row 86
$ws.Range("H86").Value = 5152.6294
$ws.Range("I86").Value = 1453.7646
$ws.Range("J86").Value = 11440.7
$ws.Range("K86").Value = 1453.7646
$ws.Range("L86").Value = 11440.7
$ws.Range("M86").Value = -330.7646
$ws.Range("N86").Value = -13686.7

# Row 89
$ws.Range("H89").Value = 5152.6294
$ws.Range("I89").Value = 1453.7646
$ws.Range("J89").Value = 11440.7
$ws.Range("K89").Value = 7268.823
$ws.Range("L89").Value = 57203.5
$ws.Range("M89").Value = -1652.823
$ws.Range("N89").Value = -68435.5

# Row 94
$ws.Range("H94").Value = 1721.8
$ws.Range("I94").Value = 1649.875
$ws.Range("J94").Value = 2009.5
$ws.Range("K94").Value = 1649.875
$ws.Range("L94").Value = 2009.5
$ws.Range("M94").Value = -1198.875
$ws.Range("N94").Value = -2911.5

# Row 105
$ws.Range("H105").Value = 3192.8215
$ws.Range("I105").Value = 2969.95
$ws.Range("J105").Value = 3750
$ws.Range("K105").Value = 2969.95
$ws.Range("L105").Value = 3750
$ws.Range("M105").Value = -1222.95
$ws.Range("N105").Value = -7244

# Row 107
$ws.Range("H107").Value = 789.6923
$ws.Range("I107").Value = 801.45
$ws.Range("J107").Value = 750.5
$ws.Range("K107").Value = 801.45
$ws.Range("L107").Value = 750.5
$ws.Range("M107").Value = 1118.55
$ws.Range("N107").Value = -4590.5

# Row 133
$ws.Range("H133").Value = 39900
$ws.Range("J133").Value = 39900
$ws.Range("L133").Value = 39900
$ws.Range("N133").Value = -50020

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1942.8572
$ws.Range("I16").Value = 666.6667
$ws.Range("J16").Value = 2900
$ws.Range("K16").Value = 666.6667
$ws.Range("L16").Value = 2900
$ws.Range("M16").Value = -379.6667
$ws.Range("N16").Value = -3474

# Row 31
$ws.Range("H31").Value = 1422
$ws.Range("I31").Value = 863.7857
$ws.Range("J31").Value = 2724.5
$ws.Range("K31").Value = 863.7857
$ws.Range("L31").Value = 2724.5
$ws.Range("M31").Value = -568.7857
$ws.Range("N31").Value = -3314.5

# Row 33
$ws.Range("H33").Value = 33010.332
$ws.Range("I33").Value = 33010.332
$ws.Range("K33").Value = 33010.332
$ws.Range("M33").Value = -32631.332

# Row 34
$ws.Range("H34").Value = 1422
$ws.Range("I34").Value = 863.7857
$ws.Range("J34").Value = 2724.5
$ws.Range("K34").Value = 863.7857
$ws.Range("L34").Value = 2724.5
$ws.Range("M34").Value = -661.7857
$ws.Range("N34").Value = -3128.5

# Row 105
$ws.Range("H105").Value = 742.6818
$ws.Range("I105").Value = 746.6111
$ws.Range("K105").Value = 746.6111
$ws.Range("M105").Value = 1000.3889

# Row 107
$ws.Range("H107").Value = 305.45456
$ws.Range("I107").Value = 234
$ws.Range("J107").Value = 346.2857
$ws.Range("K107").Value = 234
$ws.Range("L107").Value = 346.2857
$ws.Range("M107").Value = 1686
$ws.Range("N107").Value = -4186.2857

# Row 113
$ws.Range("H113").Value = 1942.8572
$ws.Range("I113").Value = 666.6667
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 666.6667
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = 1503.3333
$ws.Range("N113").Value = -7240

# Row 132
$ws.Range("H132").Value = 2657.1052
$ws.Range("I132").Value = 2218.4285
$ws.Range("J132").Value = 3885.4
$ws.Range("K132").Value = 6655.2855
$ws.Range("L132").Value = 11656.2
$ws.Range("M132").Value = -4125.2855
$ws.Range("N132").Value = -16716.2

# Row 134
$ws.Range("H134").Value = 3697.875
$ws.Range("I134").Value = 1917.6428
$ws.Range("J134").Value = 6190.2
$ws.Range("K134").Value = 5752.928400000001
$ws.Range("L134").Value = 18570.6
$ws.Range("M134").Value = -3217.928400000001
$ws.Range("N134").Value = -23640.6

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1574.421
$ws.Range("J131").Value = 1755.7551
$ws.Range("L131").Value = 5267.2653
$ws.Range("N131").Value = -15347.2653

$ws = $wb.Worksheets.Item("GSM")
# Row 140
$ws.Range("H140").Value = 57250
$ws.Range("J140").Value = 57250
$ws.Range("L140").Value = 57250
$ws.Range("N140").Value = -67610

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2995.5925
$ws.Range("I7").Value = 1966.8334
$ws.Range("J7").Value = 3289.524
$ws.Range("K7").Value = 1966.8334
$ws.Range("L7").Value = 3289.524
$ws.Range("M7").Value = -1854.8334
$ws.Range("N7").Value = -3513.524

# Row 29
$ws.Range("H29").Value = 29000
$ws.Range("I29").Value = 50000
$ws.Range("K29").Value = 50000
$ws.Range("M29").Value = -49705

# Row 40
$ws.Range("H40").Value = 3551.6296
$ws.Range("I40").Value = 1911.5555
$ws.Range("J40").Value = 4371.6665
$ws.Range("K40").Value = 1911.5555
$ws.Range("L40").Value = 4371.6665
$ws.Range("M40").Value = -1775.5555
$ws.Range("N40").Value = -4643.6665

# Row 50
$ws.Range("H50").Value = 12200
$ws.Range("I50").Value = 1000
$ws.Range("K50").Value = 1000
$ws.Range("M50").Value = -363

# Row 93
$ws.Range("H93").Value = 909.3043
$ws.Range("I93").Value = 661.875
$ws.Range("J93").Value = 1474.8572
$ws.Range("K93").Value = 661.875
$ws.Range("L93").Value = 1474.8572
$ws.Range("M93").Value = 586.125
$ws.Range("N93").Value = -3970.8572

# Row 122
$ws.Range("H122").Value = 3540.682
$ws.Range("I122").Value = 2202
$ws.Range("J122").Value = 3934.4119
$ws.Range("K122").Value = 6606
$ws.Range("L122").Value = 11803.2357
$ws.Range("M122").Value = -4156
$ws.Range("N122").Value = -16703.2357

# Row 126
$ws.Range("H126").Value = 2995.5925
$ws.Range("I126").Value = 1966.8334
$ws.Range("J126").Value = 3289.524
$ws.Range("K126").Value = 5900.5002
$ws.Range("L126").Value = 9868.572
$ws.Range("M126").Value = -3430.5002
$ws.Range("N126").Value = -14808.572

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 937.03705
$ws.Range("I107").Value = 802.9375
$ws.Range("J107").Value = 1132.091
$ws.Range("K107").Value = 2408.8125
$ws.Range("L107").Value = 3396.273
$ws.Range("M107").Value = -488.8125
$ws.Range("N107").Value = -7236.272999999999

# Row 136
$ws.Range("H136").Value = 8799558
$ws.Range("I136").Value = 9834286
$ws.Range("J136").Value = 4363.5
$ws.Range("K136").Value = 29502858
$ws.Range("L136").Value = 13090.5
$ws.Range("M136").Value = -29500308
$ws.Range("N136").Value = -18190.5
